# [Update] Modify CahierLoba following group discussion
#
# Group-discussion follow-up on the "Planification des tâches" sheet:
#  - Row 37 gets re-worded (plural "opcodes") and now has a duration.
#  - Rows 38/40/42/44/46 get their missing "Durée" values filled in.
#  - Row 39's owner changes from "Isaac" to "Tous" and gets a remark about
#    implementing/verifying the new opcodes.
#  - Row 43 gets a remark about the bulk of the work being the opcodes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39: owner -> Tous, remark added -------------------------------------
$ws.Range("C39").Value() = "Tous"
$ws.Range("E39").Value() = "Implémentation des opcodes afin de vérifier le bon fonctionnement."

# --- Row 37: new wording + 20 min duration ---------------------------------
$ws.Range("B37").Value() = "Ajouter les nouveaux opcodes dans le document sorties séquenceur"
$ws.Range("D37").Value() = 0.013888888888888888
$ws.Range("D37").NumberFormat() = "h:mm"

# --- Row 43: remark added (wrap text, top-aligned like the other remarks) ----
$ws.Range("E43").Value() = "La grande partie du travail était les opcodes que nous avons décidé d'ajouter."
$ws.Range("E43").VerticalAlignment = -4160

# --- Row 38: fill in duration (20 min) --------------------------------------
$ws.Range("D38").Value() = 0.013888888888888888
$ws.Range("D38").NumberFormat() = "h:mm"

# --- Row 40: fill in duration (1 h) ------------------------------------------
$ws.Range("D40").Value() = 0.041666666666666664
$ws.Range("D40").NumberFormat() = "h:mm"

# --- Row 42: fill in duration (45 min) ---------------------------------------
$ws.Range("D42").Value() = 0.03125
$ws.Range("D42").NumberFormat() = "h:mm"

# --- Row 44: fill in duration (45 min) ---------------------------------------
$ws.Range("D44").Value() = 0.03125

# --- Row 46: fill in duration (15 min) ----------------------------------------
$ws.Range("D46").Value() = 0.010416666666666666

# --- Update selection to match the latest edit point (F42) -------------------
$ws.Range("F42").Select()
